# Update the "5.b.1" worksheet to add the 2021 column (G) of data, matching
# the formatting of the existing 2020 column (F), and move the active
# selection to I26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the 2020 column (F4:F37) onto the new 2021 column
# (G4:G37) before filling in values, so every new cell inherits the same
# number format / font / alignment as its neighbour to the left.
$ws.Range("F4:F37").Copy()
$ws.Range("G4:G37").PasteSpecial(-4122)

# Header
$ws.Range("G4").Value = 2021

# Data rows (category/header rows - 6, 9, 12, 22, 26, 32 - are left blank,
# matching column F on those rows)
$ws.Range("G5").Value  = 92.994602261738635
$ws.Range("G7").Value  = 96.703290161846695
$ws.Range("G8").Value  = 90.681478172547202
$ws.Range("G10").Value = 94.137111592107743
$ws.Range("G11").Value = 92.036288528946841
$ws.Range("G13").Value = 93.333535676152664
$ws.Range("G14").Value = 83.868737074772881
$ws.Range("G15").Value = 96.414752760301269
$ws.Range("G16").Value = 93.03550750069985
$ws.Range("G17").Value = 92.867770679094491
$ws.Range("G18").Value = 92.877987892122846
$ws.Range("G19").Value = 96.944231580397187
$ws.Range("G20").Value = 97.652833077455327
$ws.Range("G21").Value = 95.839569929950628
$ws.Range("G23").Value = 91.964756617658111
$ws.Range("G24").Value = 95.6254397345499
$ws.Range("G25").Value = 86.706526971056761
$ws.Range("G27").Value = 68.07602651616979
$ws.Range("G28").Value = 89.482324868373453
$ws.Range("G29").Value = 92.677683811280133
$ws.Range("G30").Value = 97.302274668967783
$ws.Range("G31").Value = 98.334717525085637
$ws.Range("G33").Value = 86.414024436249534
$ws.Range("G34").Value = 93.331957186877816
$ws.Range("G35").Value = 94.055975974605076
$ws.Range("G36").Value = 95.4499076460569
$ws.Range("G37").Value = 96.978118072483895

# Move the active cell / selection (matches the author's last click before
# saving).
$ws.Range("I26").Select()
